$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataTopic")

# Insert a new column before column D (4th column) to make room for the new slot
$ws.Columns.Item(4).Insert()

# Set header for new column D1
$ws.Cells.Item(1, 4).Value = "topic_involves_anatomy"
